$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 28

# Column C ("25") must be stored as text (matching the rest of the sheet,
# which uses inline/shared strings for every cell, even numeric-looking
# ones). Temporarily mark it as Text so Excel doesn't coerce it to a
# number, then restore its formatting from a plain-text neighbor cell so
# no stray style gets left behind.
$cC = $ws.Cells.Item($row, 3)
$cC.NumberFormat = "@"
$cC.Value = "25"

$ws.Cells.Item($row, 1).Value = "2024-09-25T18:06:40Z"
$ws.Cells.Item($row, 2).Value = "temperature"
$ws.Cells.Item($row, 4).Value = "N/A"
$ws.Cells.Item($row, 5).Value = "N/A"
$ws.Cells.Item($row, 6).Value = "N/A"

# Normalize C28's style back to the sheet's default (same as the other
# cells) by copying formats from a cell that never had NumberFormat
# touched.
$ws.Cells.Item($row, 1).Copy()
$cC.PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = $false
